# Weekly fruit/vegetable price update.
# A new weekly price record is inserted as row 79 (pushing the existing
# rows 79-90 down to 80-91), for Ají "Inferno" - Primera quality,
# date serial 44449, at Terminal La Palmera de La Serena.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 79, shifting rows 79:90 down to 80:91.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly record.
$ws.Cells.Item(79, 1).Value = 8
$ws.Cells.Item(79, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(79, 3).Value = "Coquimbo"
$ws.Cells.Item(79, 4).Value = 44449
$ws.Cells.Item(79, 5).Value = 4
$ws.Cells.Item(79, 6).Value = 100112021
$ws.Cells.Item(79, 7).Value = "Ají"
$ws.Cells.Item(79, 8).Value = "Inferno"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 500
$ws.Cells.Item(79, 11).Value = 42000
$ws.Cells.Item(79, 12).Value = 43000
$ws.Cells.Item(79, 13).Value = 42500
$ws.Cells.Item(79, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 3542
$ws.Cells.Item(79, 17).Value = 12
$ws.Cells.Item(79, 18).Value = "Hortaliza"
